$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-12-09 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-10 Wednesday", 2)

# Update the division-problem table cells by explicit (row, column) position,
# since several values repeat elsewhere in the table and a plain text search
# would not target the correct occurrence.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="57÷3=19, 0"},
    @{Row=1;  Col=2; Text="24÷9=2, 6"},
    @{Row=1;  Col=3; Text="47÷9=5, 2"},
    @{Row=1;  Col=4; Text="26÷9=2, 8"},
    @{Row=1;  Col=5; Text="33÷4=8, 1"},

    @{Row=5;  Col=1; Text="15÷6=2, 3"},
    @{Row=5;  Col=2; Text="88÷8=11, 0"},
    @{Row=5;  Col=3; Text="12÷9=1, 3"},
    @{Row=5;  Col=4; Text="54÷2=27, 0"},
    @{Row=5;  Col=5; Text="54÷8=6, 6"},

    @{Row=9;  Col=1; Text="77÷5=15, 2"},
    @{Row=9;  Col=2; Text="62÷2=31, 0"},
    @{Row=9;  Col=3; Text="78÷5=15, 3"},
    @{Row=9;  Col=4; Text="64÷5=12, 4"},
    @{Row=9;  Col=5; Text="60÷7=8, 4"},

    @{Row=13; Col=1; Text="91÷5=18, 1"},
    @{Row=13; Col=2; Text="85÷8=10, 5"},
    @{Row=13; Col=3; Text="12÷7=1, 5"},
    @{Row=13; Col=4; Text="44÷2=22, 0"},
    @{Row=13; Col=5; Text="76÷3=25, 1"},

    @{Row=17; Col=1; Text="79÷2=39, 1"},
    @{Row=17; Col=2; Text="30÷5=6, 0"},
    @{Row=17; Col=3; Text="38÷5=7, 3"},
    @{Row=17; Col=4; Text="83÷9=9, 2"},
    @{Row=17; Col=5; Text="43÷4=10, 3"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
